$d = $word.ActiveDocument
$payload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16sdtdh="http://schemas.microsoft.com/office/word/2020/wordml/sdtdatahash" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex w16sdtdh wp14"><w:body><w:p w14:paraId="17BDBA0B" w14:textId="6B32EC1A" w:rsidR="0084130D" w:rsidRPr="00F12640" w:rsidRDefault="00B777AA" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00F12640"><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Resistance </w:t></w:r><w:r w:rsidR="00982C5B" w:rsidRPr="00F12640"><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>and substrate competition</w:t></w:r></w:p><w:p w14:paraId="0F8A74C8" w14:textId="4DAEC592" w:rsidR="00B777AA" w:rsidRPr="00F12640" w:rsidRDefault="00B777AA" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="44DE4343" w14:textId="1DC73749" w:rsidR="00B777AA" w:rsidRPr="00F12640" w:rsidRDefault="00B777AA" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="7EA0AC1B" w14:textId="172FD48A" w:rsidR="00B777AA" w:rsidRPr="00F12640" w:rsidRDefault="00B777AA" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00F12640"><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>PAR-3 and PKC-3</w:t></w:r></w:p><w:p w14:paraId="20F97B7D" w14:textId="1B988F11" w:rsidR="00B777AA" w:rsidRPr="00F12640" w:rsidRDefault="00B777AA" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="46B54F20" w14:textId="34A95132" w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidRDefault="00982C5B" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="5476931C" w14:textId="77777777" w:rsidR="008F0ADD" w:rsidRPr="00F12640" w:rsidRDefault="008F0ADD" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="63F686CA" w14:textId="77777777" w:rsidR="008F0ADD" w:rsidRPr="00F12640" w:rsidRDefault="008F0ADD" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00F12640"><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>PAR-2 to PAR-2</w:t></w:r></w:p><w:p w14:paraId="039DE9A0" w14:textId="218216F4" w:rsidR="008F0ADD" w:rsidRPr="00F12640" w:rsidRDefault="008F0ADD" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="102EAAD8" w14:textId="2FBD44A3" w:rsidR="008F0ADD" w:rsidRPr="00F12640" w:rsidRDefault="008F0ADD" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00F12640"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Evidence from </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>a number of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> studies suggests that, whilst PAR-2 is highly sensitive to PKC-3 when it’s weakly concentrated, it is able to resist antagonism once in an established pPAR domain. </w:t></w:r><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Although PAR-2 requires an initial aPAR asymmetry to establish a domain, stable maintenance of PAR-2 domains does not require this aPAR asymmetry to be maintained. For example, in PAR-1 knockdown/mutant conditions, aPAR and pPAR are initially segregated into domains, but </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>aPARs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> eventually return to the posterior without displacing </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>pPARs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> from the posterior cortex (Hao et al., 2006). Similarly, acute targeting of PKC-3 uniformly to the membrane is unable to fully disassemble PAR-2 domains after polarity establishment (Rodriguez et al., 2017). Interestingly, the same study showed that this resistance of PAR-2 to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>aPARs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> is acquired only when PAR-2 is concentrated in a domain (</w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>i.e.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> PAR-2 is not resistant to removal by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>aPARs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> when uniform</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> at the cortex at a lower concentration</w:t></w:r><w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">). </w:t></w:r></w:p><w:p w14:paraId="546FC58D" w14:textId="61B706E8" w:rsidR="008F0ADD" w:rsidRPr="00F12640" w:rsidRDefault="008F0ADD" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="7D47D931" w14:textId="1E0F65EC" w:rsidR="00E94F67" w:rsidRPr="00F12640" w:rsidRDefault="00E94F67" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00F12640"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">This phenomenon seems to operate in trans, rather than a cis phenomenon independent to each PAR-2 molecule. PAR-2 mutants that are usually unable to resist aPAR invasion can do so in the presence of endogenous wild-type PAR-2, indicating that PAR-2 is able to provide protection to other molecules against antagonism, possibly via a direct interaction with PKC-3, or possibly via an intermediate. The molecular details of this haven’t been determined. </w:t></w:r></w:p><w:p w14:paraId="43F48A25" w14:textId="77777777" w:rsidR="00E94F67" w:rsidRPr="00F12640" w:rsidRDefault="00E94F67" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="27DFEC3C" w14:textId="77777777" w:rsidR="008F0ADD" w:rsidRPr="00F12640" w:rsidRDefault="008F0ADD" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p w14:paraId="68C38A41" w14:textId="4B9F1333" w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidRDefault="00982C5B" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00F12640"><w:rPr><w:b/><w:bCs/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>PAR-2 to PAR-1</w:t></w:r></w:p><w:p w14:paraId="270C3B53" w14:textId="77777777" w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidRDefault="00982C5B" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="47BE3975" w14:textId="046A3136" w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidRDefault="00982C5B" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00F12640"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Part of the role of PAR-2 in recruiting PAR-1 involves interaction with PKC-3. </w:t></w:r><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">In otherwise wild type systems, PAR-1 is strictly dependent on PAR-2 to bind to the cortex, becoming entirely cytoplasmic when PAR-2 is lost (refs). In contrast, in aPAR mutant backgrounds PAR-1 shows some ability to bind to the cortex without PAR-2, but this is enhanced when PAR-2 is also present (refs). This implies a dual requirement of PAR-2 in localising PAR-1: an aPAR independent mechanism involving direct recruitment of PAR-1 by PAR-2, and a secondary role involving local protection against </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>aPARs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>.</w:t></w:r></w:p><w:p w14:paraId="67AF09EE" w14:textId="77777777" w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidRDefault="00982C5B" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="411CD3CA" w14:textId="6D342E7D" w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidRDefault="00982C5B" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>PAR-2 has been shown to inhibit phosphorylation of PAR-1 by PKC-3 in in vitro assays in a concentration-dependent manner (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Ramanujam</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">). This inhibition proceeds even in PAR-1 mutants that are unable to interact with PAR-1, implying that PAR-2 </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>is able to</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve"> act as a competitive inhibitor. That said, protection is even greater in wild type PAR-1 where an interaction with PAR-2 is permitted, suggesting that the interaction with PAR-2 can additionally block access of PKC-3 to PAR-1. </w:t></w:r><w:r w:rsidRPr="00F12640"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>It could be that PAR-2 interaction physically blocks the PKC-3 phosphorylation site on PAR-1, induces a conformational change in PAR-1 that occludes this site, or promotes membrane binding which blocks this site.</w:t></w:r></w:p><w:p w14:paraId="05E86939" w14:textId="40197444" w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidRDefault="00982C5B" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="09C46EDA" w14:textId="77777777" w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidRDefault="00982C5B" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="00F12640"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Given that competitive inhibition by PAR-2 has been described as a protective mechanism for PAR-1, it is plausible that this might contribute to maintenance of LGL-1 and CHIN-1 in a similar fashion. Indeed, whilst CHIN-1 can localise to the cortex in the absence of PAR-2 in otherwise wild type systems (unlike PAR-1), localisation is severely reduced (</w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Kumfer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F12640"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>Sailer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F12640"><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">). Alternatively, this may be a secondary consequence of rearwards flows observed in par-2 mutants (discussed </w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>earlier</w:t></w:r><w:r><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t>), or the two could be fundamentally linked.</w:t></w:r></w:p><w:p w14:paraId="4FDD51C0" w14:textId="77777777" w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidRDefault="00982C5B" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:p w14:paraId="7C4805E3" w14:textId="77777777" w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidRDefault="00982C5B" w:rsidP="00F12640"><w:pPr><w:jc w:val="both"/><w:rPr><w:color w:val="000000" w:themeColor="text1"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr></w:p><w:sectPr w:rsidR="00982C5B" w:rsidRPr="00F12640" w:rsidSect="00571055"><w:pgSz w:w="11900" w:h="16840"/><w:pgMar w:top="1440" w:right="1440" w:bottom="1440" w:left="1440" w:header="720" w:footer="720" w:gutter="0"/><w:cols w:space="720"/><w:docGrid w:linePitch="360"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Content.InsertXML($payload)
Write-Output "applied"
